# Rename the "Id" field to "SfId" and simplify its description
# from "pk, autogenerated" to "pk" on the Packet sheet (row 2: Id/SfId row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packet")

$ws.Range("F2").Value = "pk"
$ws.Range("A2").Value = "SfId"

$ws.Range("A3").Select()
